$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 816 (shifts existing 816-834 down to 820-838)
$ws.Rows("816:819").Insert()

# Row 816
$ws.Cells.Item(816, 1).Value2 = 6
$ws.Cells.Item(816, 2).Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(816, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(816, 4).Value2 = 44448
$ws.Cells.Item(816, 5).Value2 = 13
$ws.Cells.Item(816, 6).Value2 = 100112006
$ws.Cells.Item(816, 7).Value2 = 'Repollo'
$ws.Cells.Item(816, 8).Value2 = 'Crespo record'
$ws.Cells.Item(816, 9).Value2 = 'Primera'
$ws.Cells.Item(816, 10).Value2 = 13500
$ws.Cells.Item(816, 11).Value2 = 600
$ws.Cells.Item(816, 12).Value2 = 700
$ws.Cells.Item(816, 13).Value2 = 651
$ws.Cells.Item(816, 14).Value2 = '$/unidad'
$ws.Cells.Item(816, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(816, 16).Value2 = 651
$ws.Cells.Item(816, 17).Value2 = 1
$ws.Cells.Item(816, 18).Value2 = 'Hortaliza'

# Row 817
$ws.Cells.Item(817, 1).Value2 = 6
$ws.Cells.Item(817, 2).Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(817, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(817, 4).Value2 = 44448
$ws.Cells.Item(817, 5).Value2 = 13
$ws.Cells.Item(817, 6).Value2 = 100112006
$ws.Cells.Item(817, 7).Value2 = 'Repollo'
$ws.Cells.Item(817, 8).Value2 = 'Crespo record'
$ws.Cells.Item(817, 9).Value2 = 'Segunda'
$ws.Cells.Item(817, 10).Value2 = 4800
$ws.Cells.Item(817, 11).Value2 = 500
$ws.Cells.Item(817, 12).Value2 = 500
$ws.Cells.Item(817, 13).Value2 = 500
$ws.Cells.Item(817, 14).Value2 = '$/unidad'
$ws.Cells.Item(817, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(817, 16).Value2 = 500
$ws.Cells.Item(817, 17).Value2 = 1
$ws.Cells.Item(817, 18).Value2 = 'Hortaliza'

# Row 818
$ws.Cells.Item(818, 1).Value2 = 6
$ws.Cells.Item(818, 2).Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(818, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(818, 4).Value2 = 44448
$ws.Cells.Item(818, 5).Value2 = 13
$ws.Cells.Item(818, 6).Value2 = 100112006
$ws.Cells.Item(818, 7).Value2 = 'Repollo'
$ws.Cells.Item(818, 8).Value2 = 'Morada(o)'
$ws.Cells.Item(818, 9).Value2 = 'Primera'
$ws.Cells.Item(818, 10).Value2 = 4900
$ws.Cells.Item(818, 11).Value2 = 700
$ws.Cells.Item(818, 12).Value2 = 800
$ws.Cells.Item(818, 13).Value2 = 753
$ws.Cells.Item(818, 14).Value2 = '$/unidad'
$ws.Cells.Item(818, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(818, 16).Value2 = 753
$ws.Cells.Item(818, 17).Value2 = 1
$ws.Cells.Item(818, 18).Value2 = 'Hortaliza'

# Row 819
$ws.Cells.Item(819, 1).Value2 = 6
$ws.Cells.Item(819, 2).Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(819, 3).Value2 = 'Metropolitana'
$ws.Cells.Item(819, 4).Value2 = 44448
$ws.Cells.Item(819, 5).Value2 = 13
$ws.Cells.Item(819, 6).Value2 = 100112006
$ws.Cells.Item(819, 7).Value2 = 'Repollo'
$ws.Cells.Item(819, 8).Value2 = 'Morada(o)'
$ws.Cells.Item(819, 9).Value2 = 'Segunda'
$ws.Cells.Item(819, 10).Value2 = 1600
$ws.Cells.Item(819, 11).Value2 = 600
$ws.Cells.Item(819, 12).Value2 = 600
$ws.Cells.Item(819, 13).Value2 = 600
$ws.Cells.Item(819, 14).Value2 = '$/unidad'
$ws.Cells.Item(819, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(819, 16).Value2 = 600
$ws.Cells.Item(819, 17).Value2 = 1
$ws.Cells.Item(819, 18).Value2 = 'Hortaliza'

Write-Host "Edit complete."
